$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text of cell B8 to include "with the news story"
$ws.Range("B8").Value = "I also disagree with the news story. "

# Update the selection to match the target state
$ws.Range("B17").Select()
